# 10-Jul-2021 end of day update to the petty cash book.
# Sheet1 ("Buku KAS HARIAN"-style daily ledger) receives several new
# transaction rows (29-34) that were inserted above the existing "Wages
# Expense" entry for the next day, which shifts down to row 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 28: the Wages Expense debit for 9-Jul-2021 grows by an extra 260,000.
$ws.Range("D28").Formula = "=60000+260000"

# Row 29 (new): TRANSFER BCA
$ws.Range("B29").Value = "TRANSFER BCA"
$ws.Range("D29").Formula = "=4402000+6027000+800000+1510000"

# Row 30 (new): A/R
$ws.Range("B30").Value = "A/R"
$ws.Range("C30").Formula = "=6027000+6926500"

# Row 31 (new): SALES - cash/retail
$ws.Range("B31").Value = "SALES - cash/retail"
$ws.Range("C31").Formula = "=4146725+6324775-6926500"

# Row 32 (new): SOLAR - kijang D-1682-QU
$ws.Range("B32").Value = "SOLAR - kijang D-1682-QU"
$ws.Range("D32").Value = 300000

# Row 33 (new): SELISIH - lebih
$ws.Range("B33").Value = "SELISIH - lebih"
$ws.Range("C33").Value = 29500

# Row 34 (new): SETOR KE BANK
$ws.Range("B34").Value = "SETOR KE BANK"
$ws.Range("D34").Value = 3000000

# Row 35: the next day's (10-Jul-2021) Wages Expense entry, pushed down
# from row 28's neighbourhood by the six new rows above.
$ws.Range("A35").Value = 44386
$ws.Range("B35").Value = "Wages Expense"

# Update the frozen-pane scroll position / active selection to match
# where the user ended up after the day's entries (row 34 at the top of
# the scrollable area, D35 selected).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 34
$ws.Range("D35").Select()
